# Auto-generated script applying scraped-data updates to the Leve profit sheets.
$wb = $excel.ActiveWorkbook

# sheet1 (Worksheets.Item(1))
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(6, 8).Value = 625.5
$ws.Cells.Item(6, 9).Value = 667.3333
$ws.Cells.Item(6, 11).Value = 2001.9999
$ws.Cells.Item(6, 13).Value = -1889.9999
$ws.Cells.Item(62, 8).Value = 7433.1333
$ws.Cells.Item(62, 10).Value = 14449
$ws.Cells.Item(62, 12).Value = 14449
$ws.Cells.Item(62, 14).Value = -15697
$ws.Cells.Item(65, 8).Value = 7433.1333
$ws.Cells.Item(65, 10).Value = 14449
$ws.Cells.Item(65, 12).Value = 72245
$ws.Cells.Item(65, 14).Value = -78485
$ws.Cells.Item(98, 8).Value = 15248.833
$ws.Cells.Item(98, 9).Value = 13996.5
$ws.Cells.Item(98, 10).Value = 15875
$ws.Cells.Item(98, 11).Value = 13996.5
$ws.Cells.Item(98, 12).Value = 15875
$ws.Cells.Item(98, 13).Value = -12498.5
$ws.Cells.Item(98, 14).Value = -18871
$ws.Cells.Item(122, 8).Value = 15248.833
$ws.Cells.Item(122, 9).Value = 13996.5
$ws.Cells.Item(122, 10).Value = 15875
$ws.Cells.Item(122, 11).Value = 41989.5
$ws.Cells.Item(122, 12).Value = 47625
$ws.Cells.Item(122, 13).Value = -39539.5
$ws.Cells.Item(122, 14).Value = -52525
$ws.Cells.Item(141, 8).Value = 1004
$ws.Cells.Item(141, 9).Value = 1023.25
$ws.Cells.Item(141, 11).Value = 3069.75
$ws.Cells.Item(141, 13).Value = 2110.25

# sheet2 (Worksheets.Item(2))
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(32, 8).Value = 7991.6665
$ws.Cells.Item(32, 9).Value = 7991.6665
$ws.Cells.Item(32, 11).Value = 7991.6665
$ws.Cells.Item(32, 13).Value = -7704.6665
$ws.Cells.Item(63, 8).Value = 6744.3335
$ws.Cells.Item(63, 10).Value = 0
$ws.Cells.Item(63, 12).Value = 0
$ws.Cells.Item(63, 14).ClearContents()
$ws.Cells.Item(66, 8).Value = 6744.3335
$ws.Cells.Item(66, 10).Value = 0
$ws.Cells.Item(66, 12).Value = 0
$ws.Cells.Item(66, 14).ClearContents()
$ws.Cells.Item(74, 8).Value = 5624.875
$ws.Cells.Item(74, 9).Value = 4499.8
$ws.Cells.Item(74, 10).Value = 7500
$ws.Cells.Item(74, 11).Value = 4499.8
$ws.Cells.Item(74, 12).Value = 7500
$ws.Cells.Item(74, 13).Value = -3625.8
$ws.Cells.Item(74, 14).Value = -9248
$ws.Cells.Item(77, 8).Value = 5624.875
$ws.Cells.Item(77, 9).Value = 4499.8
$ws.Cells.Item(77, 10).Value = 7500
$ws.Cells.Item(77, 11).Value = 22499
$ws.Cells.Item(77, 12).Value = 37500
$ws.Cells.Item(77, 13).Value = -18131
$ws.Cells.Item(77, 14).Value = -46236

# sheet3 (Worksheets.Item(3))
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(82, 8).Value = 17851.666
$ws.Cells.Item(82, 9).Value = 6877.5
$ws.Cells.Item(82, 10).Value = 39800
$ws.Cells.Item(82, 11).Value = 6877.5
$ws.Cells.Item(82, 12).Value = 39800
$ws.Cells.Item(82, 13).Value = -6494.5
$ws.Cells.Item(82, 14).Value = -40566
$ws.Cells.Item(85, 8).Value = 17851.666
$ws.Cells.Item(85, 9).Value = 6877.5
$ws.Cells.Item(85, 10).Value = 39800
$ws.Cells.Item(85, 11).Value = 6877.5
$ws.Cells.Item(85, 12).Value = 39800
$ws.Cells.Item(85, 13).Value = -5551.5
$ws.Cells.Item(85, 14).Value = -42452
$ws.Cells.Item(86, 8).Value = 2018.6666
$ws.Cells.Item(86, 9).Value = 2018.6666
$ws.Cells.Item(86, 11).Value = 2018.6666
$ws.Cells.Item(86, 13).Value = -895.6666
$ws.Cells.Item(89, 8).Value = 2018.6666
$ws.Cells.Item(89, 9).Value = 2018.6666
$ws.Cells.Item(89, 11).Value = 10093.333
$ws.Cells.Item(89, 13).Value = -4477.333000000001

# sheet4 (Worksheets.Item(4))
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(7, 8).Value = 228.58824
$ws.Cells.Item(7, 9).Value = 262.46155
$ws.Cells.Item(7, 11).Value = 262.46155
$ws.Cells.Item(7, 13).Value = -149.46155
$ws.Cells.Item(31, 8).Value = 8702.4
$ws.Cells.Item(31, 9).Value = 3956
$ws.Cells.Item(31, 10).Value = 11075.6
$ws.Cells.Item(31, 11).Value = 3956
$ws.Cells.Item(31, 12).Value = 11075.6
$ws.Cells.Item(31, 13).Value = -3661
$ws.Cells.Item(31, 14).Value = -11665.6
$ws.Cells.Item(34, 8).Value = 8702.4
$ws.Cells.Item(34, 9).Value = 3956
$ws.Cells.Item(34, 10).Value = 11075.6
$ws.Cells.Item(34, 11).Value = 3956
$ws.Cells.Item(34, 12).Value = 11075.6
$ws.Cells.Item(34, 13).Value = -3754
$ws.Cells.Item(34, 14).Value = -11479.6
$ws.Cells.Item(58, 8).Value = 2919.4546
$ws.Cells.Item(58, 9).Value = 3086.7
$ws.Cells.Item(58, 10).Value = 1247
$ws.Cells.Item(58, 11).Value = 3086.7
$ws.Cells.Item(58, 12).Value = 1247
$ws.Cells.Item(58, 13).Value = -2883.7
$ws.Cells.Item(58, 14).Value = -1653
$ws.Cells.Item(122, 8).Value = 3999.5
$ws.Cells.Item(122, 9).Value = 3999.5
$ws.Cells.Item(122, 11).Value = 11998.5
$ws.Cells.Item(122, 13).Value = -9548.5
$ws.Cells.Item(132, 8).Value = 3040.75
$ws.Cells.Item(132, 9).Value = 3378.5386
$ws.Cells.Item(132, 10).Value = 1577
$ws.Cells.Item(132, 11).Value = 10135.6158
$ws.Cells.Item(132, 12).Value = 4731
$ws.Cells.Item(132, 13).Value = -7605.6158
$ws.Cells.Item(132, 14).Value = -9791
$ws.Cells.Item(134, 8).Value = 3055.7693
$ws.Cells.Item(134, 9).Value = 2188.6667
$ws.Cells.Item(134, 11).Value = 6566.000100000001
$ws.Cells.Item(134, 13).Value = -4031.000100000001
$ws.Cells.Item(136, 8).Value = 2919.4546
$ws.Cells.Item(136, 9).Value = 3086.7
$ws.Cells.Item(136, 10).Value = 1247
$ws.Cells.Item(136, 11).Value = 9260.099999999999
$ws.Cells.Item(136, 12).Value = 3741
$ws.Cells.Item(136, 13).Value = -6710.099999999999
$ws.Cells.Item(136, 14).Value = -8841
$ws.Cells.Item(141, 8).Value = 894812.5
$ws.Cells.Item(141, 10).Value = 894812.5
$ws.Cells.Item(141, 12).Value = 894812.5
$ws.Cells.Item(141, 14).Value = -905172.5

# sheet5 (Worksheets.Item(5))
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(122, 8).Value = 500
$ws.Cells.Item(122, 9).Value = 500
$ws.Cells.Item(122, 11).Value = 4500
$ws.Cells.Item(122, 13).Value = -2050
$ws.Cells.Item(137, 8).Value = 2000
$ws.Cells.Item(137, 9).Value = 2000
$ws.Cells.Item(137, 10).Value = 0
$ws.Cells.Item(137, 11).Value = 6000
$ws.Cells.Item(137, 12).Value = 0
$ws.Cells.Item(137, 13).Value = -900
$ws.Cells.Item(137, 14).ClearContents()

# sheet6 (Worksheets.Item(6))
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(26, 8).Value = 40000
$ws.Cells.Item(26, 9).Value = 0
$ws.Cells.Item(26, 10).Value = 40000
$ws.Cells.Item(26, 11).Value = 0
$ws.Cells.Item(26, 12).Value = 40000
$ws.Cells.Item(26, 13).ClearContents()
$ws.Cells.Item(26, 14).Value = -40560
$ws.Cells.Item(50, 8).Value = 40000
$ws.Cells.Item(50, 9).Value = 0
$ws.Cells.Item(50, 10).Value = 40000
$ws.Cells.Item(50, 11).Value = 0
$ws.Cells.Item(50, 12).Value = 40000
$ws.Cells.Item(50, 13).ClearContents()
$ws.Cells.Item(50, 14).Value = -40996
$ws.Cells.Item(80, 8).Value = 2500
$ws.Cells.Item(80, 9).Value = 2500
$ws.Cells.Item(80, 11).Value = 2500
$ws.Cells.Item(80, 13).Value = -1502
$ws.Cells.Item(83, 8).Value = 2500
$ws.Cells.Item(83, 9).Value = 2500
$ws.Cells.Item(83, 11).Value = 12500
$ws.Cells.Item(83, 13).Value = -7508
$ws.Cells.Item(97, 8).Value = 516
$ws.Cells.Item(97, 9).Value = 355
$ws.Cells.Item(97, 10).Value = 999
$ws.Cells.Item(97, 11).Value = 355
$ws.Cells.Item(97, 12).Value = 999
$ws.Cells.Item(97, 13).Value = 141
$ws.Cells.Item(97, 14).Value = -1991
$ws.Cells.Item(113, 8).Value = 3750
$ws.Cells.Item(113, 9).Value = 3750
$ws.Cells.Item(113, 11).Value = 3750
$ws.Cells.Item(113, 13).Value = -1580

# sheet8 (Worksheets.Item(8))
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(136, 8).Value = 3382.8572
$ws.Cells.Item(136, 9).Value = 3382.8572
$ws.Cells.Item(136, 11).Value = 10148.5716
$ws.Cells.Item(136, 13).Value = -7598.571599999999
